$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) / Volume(1h) (E) data block to stay plain text.
# Left to its own devices, Excel auto-coerces numeric-looking strings
# ("0.02280", "1.002", "24.830.69", "5.920") into actual numbers,
# silently dropping significant trailing zeros or mangling the value
# (these source cells use "." as a thousands separator, not a decimal
# point, e.g. "24.830.69").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.830.69"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.664.29"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "329.89"
$ws.Range("E5").Value = "  +8.38%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.3647"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("D8").Value = "47.24"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "1.136"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").Value = "0.07066"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "6.064"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "19.47"
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("D15").Value = "1.667.02"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "0.06644"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D20").Value = "78.40"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").Value = "5.920"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "15.77"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  +5.16%  "
$ws.Range("D24").Value = "24.803.96"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "2.468"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").Value = "2.427"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").Value = "148.87"
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "1.849.77"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").Value = "125.87"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "1.166"
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "5.684"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "0.08461"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").Value = "1.646"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "12.13"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "0.06247"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").Value = "0.02280"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").Value = "1.245"
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("D41").Value = "0.2084"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "8.219"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "0.5919"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "3.846"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "13.30"
$ws.Range("E46").Value = "  +6.23%  "
$ws.Range("D47").Value = "0.5657"
$ws.Range("E47").Value = "  +2.77%  "
$ws.Range("D48").Value = "125.65"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").Value = "1.945"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "0.06975"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  +4.80%  "
